$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "Maolin" run at the end of the author-list paragraph
# (the paragraph currently reads:
#   "Fangxing, Noel, Pablo, Ingo, Feili, Shunya, Maolin")
# ------------------------------------------------------------------
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("Maolin", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Re-materialize a plain Range from the Find hit's numeric bounds
    # (InsertXML must be called on a "fresh" Range, not the Find range
    # itself, otherwise content gets duplicated instead of replaced).
    $maolinRng = $d.Range($findRange.Start, $findRange.End)

    # Split the trailing ", Maolin" text so that "Maolin" becomes its own
    # spell-checked run (matching the proofErr wrapping already used for
    # the other author names such as Fangxing / Feili / Shunya in this
    # paragraph).
    $maolinXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="nb-NO"/></w:rPr><w:t>Maolin</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $maolinRng.InsertXML($maolinXml)

    # The paragraph containing the author names.
    $p2 = $d.Paragraphs(2)

    # Insert a brand-new centered paragraph right after the author names,
    # carrying the same Times New Roman / 18pt / nb-NO formatting, and set
    # its text to "Uib".
    $p2.Range.InsertParagraphAfter()
    $uibPara = $d.Paragraphs(3)
    $uibPara.Alignment = 1
    $uibPara.Range.Text = "Uib"
    $uibPara.Range.Font.Name = "Times New Roman"
    $uibPara.Range.Font.Size = 18
    $uibPara.Range.LanguageID = "nb-NO"
}
